# travail2.docx edit -- "Ajoute la generation au hazard d'un objet"
#
# The canonical-XML diff for this commit shows two related changes:
#
#   1. The <w:bookmarkStart/bookmarkEnd w:name="_GoBack"/> pair (Word's
#      "last edit position" marker) is removed from the end of the
#      "... les noms doivent s'afficher a l'ecran" paragraph.
#
#   2. That same _GoBack bookmark re-appears wrapping the paragraph
#      "generer un type d'objets au hasard : de, carte, jeton positif ou
#      jeton negatif" - i.e. this is where the author's cursor ended up
#      after editing - and that whole paragraph (the paragraph-mark run
#      properties plus all four runs) gets the same faded grey theme font
#      colour (RGB A6A6A6 / theme "background1" / shade A6) that the
#      neighbouring "commented out" checklist bullets already use.
#
# Net effect: the bullet about generating a random object type is greyed
# out like its siblings, and Word's _GoBack bookmark now marks that spot.

$d = $word.ActiveDocument

# --- Step 1: drop the stale _GoBack bookmark -------------------------------
# Bookmarks.Item resolves "_GoBack" by name even though (like real Word) it
# is hidden from Bookmarks.Count / enumeration.
$oldMark = $null
try {
    $oldMark = $d.Bookmarks.Item("_GoBack")
} catch {
    $oldMark = $null
}
if ($oldMark -ne $null) {
    $oldMark.Delete()
}

# --- Step 2: find the "generer ... au hasard ..." paragraph ---------------
$targetText = 'générer un type d’objets au hasard : dé, carte, jeton positif ou jeton négatif'
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $targetText) {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not find the 'generer un type d'objets au hasard' paragraph"
}

# --- Step 3: rewrite that paragraph (same text/structure) adding the grey --
# theme colour to every run + the paragraph mark, and re-planting _GoBack
# around the whole paragraph.
$frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00B57800" w:rsidRDefault="005150EC" w:rsidP="00875E74"><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="19"/></w:numPr><w:spacing w:line="276" w:lineRule="auto"/><w:ind w:left="851" w:right="48"/><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>générer</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> un type d’objets</w:t></w:r><w:r w:rsidR="00B57800"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> au hasard : </w:t></w:r><w:r w:rsidR="00A8698D"><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="A6A6A6" w:themeColor="background1" w:themeShade="A6"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>dé, carte, jeton positif ou jeton négatif</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.Range.InsertXML($frag)

Write-Output "Applied _GoBack relocation + grey theme colour to the 'au hasard' bullet."
